# Scheduled runner update: refresh market-board derived profit figures
# (currentAveragePrice*, LevePrice*, LeveProfit*) across several leve rows
# in multiple sheets. Mirrors an external pricing-tool data refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("N7").ClearContents()

$ws.Range("H14").Value = 0
$ws.Range("J14").Value = 0
$ws.Range("L14").Value = 0
$ws.Range("N14").ClearContents()

$ws.Range("H74").Value = 4433.533
$ws.Range("I74").Value = 4386.143
$ws.Range("J74").Value = 4475
$ws.Range("K74").Value = 4386.143
$ws.Range("L74").Value = 4475
$ws.Range("M74").Value = -3450.143
$ws.Range("N74").Value = -6347

$ws.Range("H77").Value = 4433.533
$ws.Range("I77").Value = 4386.143
$ws.Range("J77").Value = 4475
$ws.Range("K77").Value = 21930.715
$ws.Range("L77").Value = 22375
$ws.Range("M77").Value = -17250.715
$ws.Range("N77").Value = -31735

$ws.Range("H132").Value = 52836.9
$ws.Range("I132").Value = 55449.42
$ws.Range("J132").Value = 3199
$ws.Range("K132").Value = 166348.26
$ws.Range("L132").Value = 9597
$ws.Range("M132").Value = -163818.26
$ws.Range("N132").Value = -14657

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 7079.636
$ws.Range("I88").Value = 11215.272
$ws.Range("J88").Value = 2944
$ws.Range("K88").Value = 11215.272
$ws.Range("L88").Value = 2944
$ws.Range("M88").Value = -10809.272
$ws.Range("N88").Value = -3756

$ws.Range("H91").Value = 7079.636
$ws.Range("I91").Value = 11215.272
$ws.Range("J91").Value = 2944
$ws.Range("K91").Value = 11215.272
$ws.Range("L91").Value = 2944
$ws.Range("M91").Value = -9811.272000000001
$ws.Range("N91").Value = -5752

$ws.Range("H112").Value = 19512.857
$ws.Range("J112").Value = 19512.857
$ws.Range("L112").Value = 19512.857
$ws.Range("N112").Value = -22466.857

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H100").Value = 24481
$ws.Range("I100").Value = 19800
$ws.Range("J100").Value = 26821.5
$ws.Range("K100").Value = 19800
$ws.Range("L100").Value = 26821.5
$ws.Range("M100").Value = -18718
$ws.Range("N100").Value = -28985.5

$ws.Range("H141").Value = 44612.5
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 44612.5
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 44612.5
$ws.Range("N141").Value = -54972.5
$ws.Range("M141").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3553.1807
$ws.Range("I31").Value = 910.10345
$ws.Range("J31").Value = 4972.6113
$ws.Range("K31").Value = 910.10345
$ws.Range("L31").Value = 4972.6113
$ws.Range("M31").Value = -615.10345
$ws.Range("N31").Value = -5562.6113

$ws.Range("H34").Value = 3553.1807
$ws.Range("I34").Value = 910.10345
$ws.Range("J34").Value = 4972.6113
$ws.Range("K34").Value = 910.10345
$ws.Range("L34").Value = 4972.6113
$ws.Range("M34").Value = -708.10345
$ws.Range("N34").Value = -5376.6113

$ws.Range("H62").Value = 5602.643
$ws.Range("I62").Value = 5344.2856
$ws.Range("J62").Value = 5861
$ws.Range("K62").Value = 5344.2856
$ws.Range("L62").Value = 5861
$ws.Range("M62").Value = -4720.2856
$ws.Range("N62").Value = -7109

$ws.Range("H65").Value = 5602.643
$ws.Range("I65").Value = 5344.2856
$ws.Range("J65").Value = 5861
$ws.Range("K65").Value = 26721.428
$ws.Range("L65").Value = 29305
$ws.Range("M65").Value = -23601.428
$ws.Range("N65").Value = -35545

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 651509.3
$ws.Range("I5").Value = 1029.25
$ws.Range("J5").Value = 764636.3
$ws.Range("K5").Value = 3087.75
$ws.Range("L5").Value = 2293908.9
$ws.Range("M5").Value = -2975.75
$ws.Range("N5").Value = -2294132.9

$ws.Range("H9").Value = 207999.8
$ws.Range("J9").Value = 207999.8
$ws.Range("L9").Value = 623999.3999999999
$ws.Range("N9").Value = -624447.3999999999

$ws.Range("H68").Value = 976.23
$ws.Range("I68").Value = 673.65
$ws.Range("J68").Value = 1430.1
$ws.Range("K68").Value = 2020.95
$ws.Range("L68").Value = 4290.299999999999
$ws.Range("M68").Value = -1209.95
$ws.Range("N68").Value = -5912.299999999999

$ws.Range("H71").Value = 976.23
$ws.Range("I71").Value = 673.65
$ws.Range("J71").Value = 1430.1
$ws.Range("K71").Value = 6062.849999999999
$ws.Range("L71").Value = 12870.9
$ws.Range("M71").Value = -2006.849999999999
$ws.Range("N71").Value = -20982.9

$ws.Range("H98").Value = 566
$ws.Range("I98").Value = 615.5
$ws.Range("K98").Value = 1846.5
$ws.Range("M98").Value = -348.5

$ws.Range("H131").Value = 3612.7234
$ws.Range("I131").Value = 662.375
$ws.Range("J131").Value = 4217.923
$ws.Range("K131").Value = 1987.125
$ws.Range("L131").Value = 12653.769
$ws.Range("M131").Value = 3052.875
$ws.Range("N131").Value = -22733.769

$ws.Range("H132").Value = 1344.862
$ws.Range("I132").Value = 1446.5714
$ws.Range("J132").Value = 1249.9333
$ws.Range("K132").Value = 13019.1426
$ws.Range("L132").Value = 11249.3997
$ws.Range("M132").Value = -10489.1426
$ws.Range("N132").Value = -16309.3997

$ws.Range("H135").Value = 651509.3
$ws.Range("I135").Value = 1029.25
$ws.Range("J135").Value = 764636.3
$ws.Range("K135").Value = 9263.25
$ws.Range("L135").Value = 6881726.7
$ws.Range("M135").Value = -6728.25
$ws.Range("N135").Value = -6886796.7

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()

$ws.Range("H63").Value = 5555
$ws.Range("J63").Value = 5555
$ws.Range("L63").Value = 5555
$ws.Range("N63").Value = -6927

$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()

$ws.Range("H66").Value = 5555
$ws.Range("J66").Value = 5555
$ws.Range("L66").Value = 16665
$ws.Range("N66").Value = -23529

$ws.Range("H132").Value = 4682.8237
$ws.Range("I132").Value = 4445.5
$ws.Range("J132").Value = 4755.846
$ws.Range("K132").Value = 13336.5
$ws.Range("L132").Value = 14267.538
$ws.Range("M132").Value = -10806.5
$ws.Range("N132").Value = -19327.538

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H104").Value = 17592.5
$ws.Range("J104").Value = 16790
$ws.Range("L104").Value = 16790
$ws.Range("N104").Value = -23778

$ws.Range("H139").Value = 47473.332
$ws.Range("J139").Value = 47473.332
$ws.Range("L139").Value = 47473.332
$ws.Range("N139").Value = -57753.332

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H103").Value = 22499.5
$ws.Range("J103").Value = 22499.5
$ws.Range("L103").Value = 22499.5
$ws.Range("N103").Value = -24843.5
